# "flexen met mn punten" - update progress on the "4c" assignment (row 13)
# and the "totaal" row (row 16): move from "afwachten"/"mee bezig" status
# text to actual scored points, and flag the new totaal value with a
# conditional format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("1a"): was marked "afwachten" -> now scored 5 points.
$ws.Range("E5").Value = 5

# Row 13 ("4c"): was still open -> now marked "mee bezig".
$ws.Range("E13").Value = "mee bezig"

# Row 16 ("totaal" sub-row): was marked "mee bezig" -> now scored 10 points.
$ws.Range("E16").Value = 10

# Highlight E16 (green) when its value is greater than 1, same "Highlight
# Cell Rules > Greater Than..." style already used elsewhere on the sheet,
# placed at the top of the conditional-formatting priority stack.
$fc = $ws.Range("E16").FormatConditions.Add(1, 5, "1")
$fc.Font.Color = 24832
$fc.Interior.Color = 13561798
$fc.SetFirstPriority()

# Selection was left on I15.
$ws.Range("I15").Select()
